$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 23.331730174041
$ws.Range("B2").Value = 19.69795564090124
$ws.Range("C2").Value = 27.10557193995341
$ws.Range("A3").Value = 19.09035096275557
$ws.Range("B3").Value = 11.89577523062924
$ws.Range("C3").Value = 27.47572503725279
$ws.Range("A4").Value = 20.07239925893111
$ws.Range("B4").Value = 14.79264873210985
$ws.Range("C4").Value = 25.496635695391
$ws.Range("A5").Value = 23.4104594049059
$ws.Range("B5").Value = 16.30587851826161
$ws.Range("C5").Value = 31.10660148652546
$ws.Range("A6").Value = 28.01060954458787
$ws.Range("B6").Value = 24.11831176849653
$ws.Range("C6").Value = 32.15271988984622
$ws.Range("A7").Value = 20.80394973453817
$ws.Range("B7").Value = 18.23968215365121
$ws.Range("C7").Value = 23.79748862334326
$ws.Range("A8").Value = 10.36346117576603
$ws.Range("B8").Value = 8.076723119444825
$ws.Range("C8").Value = 13.50103445939163
$ws.Range("A9").Value = 25.92323967699945
$ws.Range("B9").Value = 22.67097920063459
$ws.Range("C9").Value = 29.50983921648628
$ws.Range("A10").Value = 10.87612183307555
$ws.Range("B10").Value = 8.342057663979126
$ws.Range("C10").Value = 14.51922200463839
$ws.Range("A11").Value = 18.62933805463966
$ws.Range("B11").Value = 14.60070701321043
$ws.Range("C11").Value = 22.75305904684553
$ws.Range("A12").Value = 30.64849035416089
$ws.Range("B12").Value = 27.30595746250612
$ws.Range("C12").Value = 34.66801265088516
$ws.Range("A13").Value = 12.83439254317619
$ws.Range("B13").Value = 9.610132574885627
$ws.Range("C13").Value = 17.68498909169671
$ws.Range("A14").Value = 14.98044654714227
$ws.Range("B14").Value = 11.36310349756332
$ws.Range("C14").Value = 19.31660022161913
$ws.Range("A15").Value = 27.76868222476446
$ws.Range("B15").Value = 19.66382510759768
$ws.Range("C15").Value = 36.79170090135661
$ws.Range("A16").Value = 27.99981108743031
$ws.Range("B16").Value = 24.12276784547966
$ws.Range("C16").Value = 32.13037533926538
$ws.Range("A17").Value = 30.57939235983462
$ws.Range("B17").Value = 27.19472622032663
$ws.Range("C17").Value = 34.29945531469322
$ws.Range("A18").Value = 20.06964453289743
$ws.Range("B18").Value = 14.79750830167271
$ws.Range("C18").Value = 25.45604483392655
$ws.Range("A19").Value = 19.64096454265627
$ws.Range("B19").Value = 15.19829840041572
$ws.Range("C19").Value = 25.39382700456649
$ws.Range("A20").Value = 32.1011149354301
$ws.Range("B20").Value = 27.69734033597354
$ws.Range("C20").Value = 36.12050282748923
$ws.Range("A21").Value = 29.01838726366836
$ws.Range("B21").Value = 22.6530271988896
$ws.Range("C21").Value = 35.60990629179418
$ws.Range("A22").Value = 11.16923821505657
$ws.Range("B22").Value = 8.354055180567201
$ws.Range("C22").Value = 14.9053032184349
$ws.Range("A23").Value = 29.35753887327363
$ws.Range("B23").Value = 21.82922772245666
$ws.Range("C23").Value = 38.57332591782405
$ws.Range("A24").Value = 27.9999640378042
$ws.Range("B24").Value = 24.12758284979184
$ws.Range("C24").Value = 32.13889450137722
$ws.Range("A25").Value = 23.53221420870001
$ws.Range("B25").Value = 19.96393819573313
$ws.Range("C25").Value = 27.23767151544986
$ws.Range("A26").Value = 29.14708631828314
$ws.Range("B26").Value = 23.64124173484163
$ws.Range("C26").Value = 35.92087925821886
$ws.Range("A27").Value = 13.34608248747216
$ws.Range("B27").Value = 9.59038717704694
$ws.Range("C27").Value = 18.55369101761913
$ws.Range("A28").Value = 22.18725896429212
$ws.Range("B28").Value = 19.26562367104054
$ws.Range("C28").Value = 25.66360521339153
$ws.Range("A29").Value = 30.55309994573313
$ws.Range("B29").Value = 22.42835031634755
$ws.Range("C29").Value = 40.01529550083274
$ws.Range("A30").Value = 28.68949159603353
$ws.Range("B30").Value = 25.46526568297984
$ws.Range("C30").Value = 32.21426435698041
$ws.Range("A31").Value = 14.46468687374846
$ws.Range("B31").Value = 11.19214214041788
$ws.Range("C31").Value = 18.41367262197492
$ws.Range("A32").Value = 11.8958219074782
$ws.Range("B32").Value = 9.170713061542362
$ws.Range("C32").Value = 15.56046125017654
$ws.Range("A33").Value = 14.97795683639712
$ws.Range("B33").Value = 11.48193882820627
$ws.Range("C33").Value = 19.37914549553492
$ws.Range("A34").Value = 26.83426623911475
$ws.Range("B34").Value = 21.95319444783581
$ws.Range("C34").Value = 32.08434856313949
$ws.Range("A35").Value = 30.04835336739003
$ws.Range("B35").Value = 25.91168828310794
$ws.Range("C35").Value = 33.99015302758899
$ws.Range("A36").Value = 20.7603942618903
$ws.Range("B36").Value = 18.17282895325192
$ws.Range("C36").Value = 23.61627565627834
$ws.Range("A37").Value = 31.46673017730801
$ws.Range("B37").Value = 25.15581067705905
$ws.Range("C37").Value = 37.66974328913387
$ws.Range("A38").Value = 15.67166887978145
$ws.Range("B38").Value = 11.80516932687377
$ws.Range("C38").Value = 21.08556443863539
$ws.Range("A39").Value = 31.92965902057797
$ws.Range("B39").Value = 28.4956885176541
$ws.Range("C39").Value = 36.04730120384583
$ws.Range("A40").Value = 24.46927609832191
$ws.Range("B40").Value = 18.24452766784713
$ws.Range("C40").Value = 30.20113101641518
$ws.Range("A41").Value = 19.31537950988023
$ws.Range("B41").Value = 16.94033284516921
$ws.Range("C41").Value = 21.8549626258784
$ws.Range("A42").Value = 19.16321496394614
$ws.Range("B42").Value = 15.03954922391687
$ws.Range("C42").Value = 24.09326640239387
$ws.Range("A43").Value = 28.70884668007166
$ws.Range("B43").Value = 25.46808689508625
$ws.Range("C43").Value = 32.25410159210425
$ws.Range("A44").Value = 16.130458118752
$ws.Range("B44").Value = 13.49653805135677
$ws.Range("C44").Value = 19.29311192974283
$ws.Range("A45").Value = 20.10379110110047
$ws.Range("B45").Value = 15.48587735590412
$ws.Range("C45").Value = 25.82402505092682
$ws.Range("A46").Value = 13.98828085632986
$ws.Range("B46").Value = 11.10194470332602
$ws.Range("C46").Value = 17.45433776624914
$ws.Range("A47").Value = 11.44448969684385
$ws.Range("B47").Value = 8.322700831999505
$ws.Range("C47").Value = 15.716800806695
$ws.Range("A48").Value = 30.589620843428
$ws.Range("B48").Value = 27.21309293626852
$ws.Range("C48").Value = 34.5676016384662
$ws.Range("A49").Value = 22.43509686776643
$ws.Range("B49").Value = 17.18808548051937
$ws.Range("C49").Value = 28.27525548320042
$ws.Range("A50").Value = 31.41082620535865
$ws.Range("B50").Value = 25.040356390542
$ws.Range("C50").Value = 37.7255636678384
$ws.Range("A51").Value = 19.00326048169884
$ws.Range("B51").Value = 11.78705100312315
$ws.Range("C51").Value = 27.38008492419764
$ws.Range("A52").Value = 13.27428428927199
$ws.Range("B52").Value = 10.68367372365218
$ws.Range("C52").Value = 16.37041806941615
$ws.Range("A53").Value = 31.79143274325667
$ws.Range("B53").Value = 28.35170953365126
$ws.Range("C53").Value = 35.83479056561359
$ws.Range("A54").Value = 27.73767226895392
$ws.Range("B54").Value = 19.61006759503502
$ws.Range("C54").Value = 36.79254696320983
$ws.Range("A55").Value = 28.95715691085401
$ws.Range("B55").Value = 22.67971830459702
$ws.Range("C55").Value = 35.48520453683066
$ws.Range("A56").Value = 14.45759393959737
$ws.Range("B56").Value = 11.2007105501913
$ws.Range("C56").Value = 18.36482559777452
$ws.Range("A57").Value = 11.8958219074782
$ws.Range("B57").Value = 9.170713061542362
$ws.Range("C57").Value = 15.56046125017654
$ws.Range("A58").Value = 31.43014321933036
$ws.Range("B58").Value = 25.46212102771626
$ws.Range("C58").Value = 37.90037748837597
$ws.Range("A59").Value = 17.92165237797689
$ws.Range("B59").Value = 14.50031589412157
$ws.Range("C59").Value = 22.13219804374561
$ws.Range("A60").Value = 29.8114479868593
$ws.Range("B60").Value = 21.80754196971159
$ws.Range("C60").Value = 38.90242099318785
$ws.Range("A61").Value = 15.10797538264709
$ws.Range("B61").Value = 11.57533191122437
$ws.Range("C61").Value = 19.90424456360322
$ws.Range("A62").Value = 18.64533292504597
$ws.Range("B62").Value = 14.61556212902492
$ws.Range("C62").Value = 22.77018461058528
$ws.Range("A63").Value = 15.80767971967518
$ws.Range("B63").Value = 11.71119395714268
$ws.Range("C63").Value = 21.18988221180542
$ws.Range("A64").Value = 29.73690682918778
$ws.Range("B64").Value = 21.80445966499818
$ws.Range("C64").Value = 38.98847222046594
$ws.Range("A65").Value = 11.44041430655556
$ws.Range("B65").Value = 8.9180167762152
$ws.Range("C65").Value = 14.30330386403321
$ws.Range("A66").Value = 25.48822261905331
$ws.Range("B66").Value = 21.1597269986093
$ws.Range("C66").Value = 30.57892666049671
$ws.Range("A67").Value = 24.49821204508587
$ws.Range("B67").Value = 20.66559853120764
$ws.Range("C67").Value = 28.50937072781771
$ws.Range("A68").Value = 26.94463322395677
$ws.Range("B68").Value = 21.23946488654743
$ws.Range("C68").Value = 32.89961661937461
$ws.Range("A69").Value = 11.44316534284396
$ws.Range("B69").Value = 8.543769693083632
$ws.Range("C69").Value = 15.70122042279455
$ws.Range("A70").Value = 24.42920459818092
$ws.Range("B70").Value = 18.5768688971132
$ws.Range("C70").Value = 30.93718498362488
$ws.Range("A71").Value = 21.16171899170495
$ws.Range("B71").Value = 18.86586499709033
$ws.Range("C71").Value = 23.85741424371978
$ws.Range("A72").Value = 25.049243825261
$ws.Range("B72").Value = 22.24693996848086
$ws.Range("C72").Value = 28.13935511120514
